$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title (row 1) with the new date
$ws.Range("A1").Value = "An Individual Customer's usage on 2025-03-03"

# Update the demand values for the remaining rows
$ws.Range("B3").Value = 0.014
$ws.Range("B4").Value = 0.007

# Remove rows 5 through 12 (times 02:00:00 .. 09:00:00) entirely,
# shrinking the used range down to A1:C4
$ws.Range("A5:C12").ClearContents()
